# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" and refresh the
# related handoff timestamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns for zh-cn (B2) and de-de (C2), plus the
# overall Latest Handoff Date (D2).
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-03-21 08:32:27"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (E2).
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-21 08:32:23"

# de-de sheet: Status (C2) and Latest Handoff Datetime (E2).
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-21 08:32:27"
